$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = 28
$ws.Range("C28").Value = 15
$ws.Range("D28").Value = 16
$ws.Range("E28").Value = 31
$ws.Range("F28").Value = 59
$ws.Range("G28").Value = 90
